# lista_espera_coach_imagen.xlsx — add the new waitlist signup (row 3) that
# came in on 4/6/2025, and make sure the sheet view is explicitly
# left-to-right (matches the "Home" rebuild mentioned in the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Luciano Albani"
$ws.Range("B3").Value = "luchoalbanix1@gmail.com"
$ws.Range("C3").Value = "+54 2234480301"
$ws.Range("D3").Value = "4/6/2025, 9:32:56 p.m."

# Explicitly pin the sheet/window to left-to-right reading order.
$ws.DisplayRightToLeft = $false
$excel.ActiveWindow.DisplayRightToLeft = $false
